# Update the 16_cues word/image/category table to the new word/image set,
# keeping the category in sync with each image (house.* -> "house", flower.* -> "flower").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 'tollen'
$ws.Range("B3").Value = 'flower/flower012.jpg'
$ws.Range("C3").Value = 'flower'

$ws.Range("A4").Value = 'platzen'
$ws.Range("B4").Value = 'flower/flower027.jpg'
$ws.Range("C4").Value = 'flower'

$ws.Range("A6").Value = 'beißen'
$ws.Range("B6").Value = 'flower/flower031.jpg'
$ws.Range("C6").Value = 'flower'

$ws.Range("A7").Value = 'kehren'
$ws.Range("B7").Value = 'house/house006.jpg'
$ws.Range("C7").Value = 'house'

$ws.Range("A9").Value = 'ändern'
$ws.Range("B9").Value = 'house/house020.jpg'
$ws.Range("C9").Value = 'house'

$ws.Range("A10").Value = 'wohnen'
$ws.Range("B10").Value = 'house/house022.jpg'
$ws.Range("C10").Value = 'house'

$ws.Range("A12").Value = 'altern'
$ws.Range("B12").Value = 'flower/flower028.jpg'
$ws.Range("C12").Value = 'flower'

$ws.Range("A13").Value = 'antun'
$ws.Range("B13").Value = 'house/house003.jpg'
$ws.Range("C13").Value = 'house'

$ws.Range("A15").Value = 'segnen'
$ws.Range("B15").Value = 'flower/flower008.jpg'
$ws.Range("C15").Value = 'flower'

$ws.Range("A16").Value = 'bellen'
$ws.Range("B16").Value = 'flower/flower032.jpg'
$ws.Range("C16").Value = 'flower'

$ws.Range("A18").Value = 'biegen'
$ws.Range("B18").Value = 'house/house016.jpg'
$ws.Range("C18").Value = 'house'

$ws.Range("A19").Value = 'gleichen'
$ws.Range("B19").Value = 'flower/flower001.jpg'
$ws.Range("C19").Value = 'flower'

$ws.Range("A21").Value = 'schicken'
$ws.Range("B21").Value = 'house/house019.jpg'
$ws.Range("C21").Value = 'house'

$ws.Range("A22").Value = 'lehnen'
$ws.Range("B22").Value = 'house/house029.jpg'
$ws.Range("C22").Value = 'house'

$ws.Range("A24").Value = 'schützen'
$ws.Range("B24").Value = 'flower/flower029.jpg'
$ws.Range("C24").Value = 'flower'

$ws.Range("A25").Value = 'süßen'
$ws.Range("B25").Value = 'house/house027.jpg'
$ws.Range("C25").Value = 'house'

$ws.Range("A27").Value = 'sammeln'
$ws.Range("B27").Value = 'house/house007.jpg'
$ws.Range("C27").Value = 'house'

$ws.Range("A28").Value = 'grüßen'
$ws.Range("B28").Value = 'flower/flower024.jpg'
$ws.Range("C28").Value = 'flower'

$ws.Range("A30").Value = 'töten'
$ws.Range("B30").Value = 'house/house009.jpg'
$ws.Range("C30").Value = 'house'

$ws.Range("A31").Value = 'rasen'
$ws.Range("B31").Value = 'flower/flower033.jpg'
$ws.Range("C31").Value = 'flower'

$ws.Range("A33").Value = 'wachsen'
$ws.Range("B33").Value = 'flower/flower003.jpg'
$ws.Range("C33").Value = 'flower'

$ws.Range("A34").Value = 'schwimmen'
$ws.Range("B34").Value = 'house/house025.jpg'
$ws.Range("C34").Value = 'house'

$ws.Range("A36").Value = 'landen'
$ws.Range("B36").Value = 'flower/flower020.jpg'
$ws.Range("C36").Value = 'flower'

$ws.Range("A37").Value = 'steuern'
$ws.Range("B37").Value = 'flower/flower014.jpg'
$ws.Range("C37").Value = 'flower'

$ws.Range("A39").Value = 'öffnen'
$ws.Range("B39").Value = 'house/house030.jpg'
$ws.Range("C39").Value = 'house'

$ws.Range("A40").Value = 'dürfen'
$ws.Range("B40").Value = 'flower/flower002.jpg'
$ws.Range("C40").Value = 'flower'

$ws.Range("A42").Value = 'werden'
$ws.Range("B42").Value = 'flower/flower006.jpg'
$ws.Range("C42").Value = 'flower'

$ws.Range("A43").Value = 'zahlen'
$ws.Range("B43").Value = 'flower/flower010.jpg'
$ws.Range("C43").Value = 'flower'

$ws.Range("A45").Value = 'stopfen'
$ws.Range("B45").Value = 'house/house011.jpg'
$ws.Range("C45").Value = 'house'

$ws.Range("A46").Value = 'bieten'
$ws.Range("B46").Value = 'house/house001.jpg'
$ws.Range("C46").Value = 'house'

$ws.Range("A48").Value = 'fragen'
$ws.Range("B48").Value = 'house/house010.jpg'
$ws.Range("C48").Value = 'house'

$ws.Range("A49").Value = 'bilden'
$ws.Range("B49").Value = 'house/house013.jpg'
$ws.Range("C49").Value = 'house'
